$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (column C) date from 45189 to 45190 for all
#    existing data rows (rows 2 through 115).
for ($r = 2; $r -le 115; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}

# 2. Row 115 gains an explicit row height (15pt / custom height), matching
#    the rest of the data rows.
$ws.Rows.Item(115).RowHeight = 15

# 3. Append a brand new row 116 with a new logging notification record.
$ws.Cells.Item(116, 1).Value = "A 44448-2023"

$ws.Cells.Item(116, 2).Value = 45189
$ws.Cells.Item(116, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(116, 3).Value = 45190
$ws.Cells.Item(116, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(116, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(116, 5).Value = "ÄNGELHOLM"

$ws.Cells.Item(116, 7).Value = 0.6
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 0
$ws.Cells.Item(116, 14).Value = 0
$ws.Cells.Item(116, 15).Value = 0
$ws.Cells.Item(116, 16).Value = 0
$ws.Cells.Item(116, 17).Value = 0

# Column R keeps the same "wrap text" style used throughout the sheet, but
# stays empty for the new row, just like every other row.
$ws.Cells.Item(115, 18).Copy($ws.Cells.Item(116, 18))
$ws.Cells.Item(116, 18).Value = ""
